# Split the old "Terms Typically Offered" column (D) apart from three new
# requirement columns (Corequisites, Concurrent, Recommended) that were
# inserted ahead of it, per "Updated data to reflect new requirement
# separation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank columns at D (old D "Terms Typically Offered"
# shifts right to G; dimension grows from A1:D22 to A1:G22).
$ws.Columns("D:F").Insert()

# New header row (row 1).
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# For every data row, the three new columns default to "NA" ...
for ($r = 2; $r -le 22; $r++) {
    $ws.Range("D${r}:F${r}").Value = "NA"
}

# ... except row 16, whose old D cell actually held a mangled value
# ("SPCorequisite: BUS 417, ...") mixing the term code with a corequisite
# note. Split it back apart: the corequisite text becomes the new D16,
# and the term code ("SP") goes into the shifted G16.
$nbsp = [char]0x00A0
$ws.Range("D16").Value = "BUS${nbsp}417, or BUS 414 and BUS 415, and OCOB graduate standing or approval from the Associate Dean."
$ws.Range("G16").Value = "SP"

# Minor data cleanup: normalize the stray non-breaking space in the GSA 543
# prerequisite text (row 9, column C) to a regular space.
$ws.Range("C9").Value = "GSA 541 and OCOB graduate standing or approval from the Associate Dean."
